$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Set the new task and estimate values
$ws.Range("A2").Value = "Specular lighting"
$ws.Range("B2").Value = 4

# The inserted row inherited the bold header style; make it explicitly
# non-bold so it gets its own (unbolded) style definition.
$ws.Range("A2:B2").Font.Bold = $false

# Select the new row
$ws.Range("A2:B2").Select()
